# Screening_overview.xlsx edit:
# Insert a new "Analyse2" column (W) so conditions can be selected more
# precisely for plotting. Existing columns W:Y (Tween_bad, DMSO_bad,
# RASTRIC_SD) shift right to X:Z. Also clears the "Analyse" flag (V) for a
# handful of rows that are now covered by the new "Analyse2" flag instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at W - this shifts the old W,X,Y (Tween_bad,
#    DMSO_bad, RASTRIC_SD) one position to the right (X,Y,Z).
$ws.Columns("W").Insert()

# 2. Header for the freshly inserted column.
$ws.Range("W1").Value = "Analyse2"

# 3. Fill the new column's data rows (2-55) with 0, then flag the two rows
#    (39 and 40 - the OPT0024 WNT_high / WNT_low conditions) that should be
#    picked out by the new "Analyse2" selector.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 23).Value = 0
}
$ws.Cells.Item(39, 23).Value = 1
$ws.Cells.Item(40, 23).Value = 1

# 4. These rows are no longer flagged under the original "Analyse" (V)
#    column - they are now represented via "Analyse2" instead.
$ws.Range("V35").Value = 0
$ws.Range("V36").Value = 0
$ws.Range("V37").Value = 0
$ws.Range("V38").Value = 0
$ws.Range("V41").Value = 0

# 5. Re-establish the AutoFilter over the now-wider table (A:Y instead of
#    A:X) and fix up the hidden _FilterDatabase defined name to match.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:Y55").AutoFilter()

foreach ($n in $wb.Names) {
    $nm = $n.Name()
    if ($nm -eq "Sheet 1!_FilterDatabase") {
        $n.RefersTo = "='Sheet 1'!`$A`$1:`$Y`$55"
    }
}

# 6. Restore view state: zoom back to 100%, and selection in the frozen
#    bottom-right pane moved from V35 to V45.
$excel.ActiveWindow.Zoom = 100
$ws.Range("V45").Select()
